$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update version strings for CQLive-PENGUIN1SA4062 from 1.23 -> 1.24
$ws.Range("I2").Value = "CQLive-PENGUIN1SA4062-1.24"
$ws.Range("J2").Value = "CQLive-PENGUIN1SA4062-1.24.apk"

# Update the active selection to J10
$ws.Range("J10").Select()
